# Atualização automática da planilha
# Updates the "Organograma" sheet: clears the "Key User BPM / Key User
# Riscos / Key User Auditorias" rows (25-29) down to just the Área value
# ("SoftExpert"), removes the "Planejamento Estratégico" área on row 30,
# bumps those rows to a 15pt height, and moves the sheet's frozen-pane /
# selection to where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organograma")

# Rows 25-29: clear Papel / Nome / Cargo (B:D), keep "Time Projeto" in A
# and set Área (E) to "SoftExpert".
$ws.Range("B25:D29").ClearContents()
$ws.Range("E25").Value = "SoftExpert"
$ws.Range("E26").Value = "SoftExpert"
$ws.Range("E27").Value = "SoftExpert"
$ws.Range("E28").Value = "SoftExpert"
$ws.Range("E29").Value = "SoftExpert"

# Row 30: Área cleared.
$ws.Range("E30").ClearContents()

# Rows 25-30 now use a 15pt row height.
$ws.Range("A25:A30").RowHeight = 15

# Move the active selection to reflect where the user left off scrolling
# (the header rows stay frozen via the sheet's existing ySplit=2 pane;
# only the in-view selection/scroll position moves).
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 20
$ws.Range("E44").Select() | Out-Null
